$wb = $excel.ActiveWorkbook

# Rename the Taxon sheet's Id value from "mpn_m129" to "taxon"
$taxonSheet = $wb.Worksheets.Item("Taxon")
$taxonSheet.Range("B1").Value = "taxon"

# Insert a new "Environment" worksheet right after "Taxon" (before "Submodels")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $taxonSheet)
$newSheet.Name = "Environment"

$newSheet.Range("A1").Value = "Id"
$newSheet.Range("B1").Value = "env"
$newSheet.Range("A2").Value = "Name"
$newSheet.Range("A3").Value = "Temperature"
$newSheet.Range("B3").Value = 37
$newSheet.Range("A4").Value = "Temperature units"
$newSheet.Range("B4").Value = "C"
$newSheet.Range("A5").Value = "pH"
$newSheet.Range("B5").Value = 7.75
$newSheet.Range("A6").Value = "Database references"
$newSheet.Range("A7").Value = "Comments"
$newSheet.Range("A8").Value = "References"

$newSheet.Select()
